# Full reworking of all models
# - Rename "table" sheet to "All_model_short_table"
# - Update coefficient/p-value data on both sheets (re-run of the model)
# - Re-shape the "table" sheet so it mirrors the "All_model_short" sheet
#   (plain A:C layout) instead of the old B:C + decorative I:J columns
# - Remove the conditional formatting (highlight p in [0.05, 0.1]) from
#   "All_model_short"
# - Tidy up now-unused custom column widths

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: parse scientific-notation strings into doubles (the PS parser
# here chokes on bare "1E-5" literals), and apply one row of (param,
# coefficient, p) data onto a worksheet starting at a given row.
# ---------------------------------------------------------------------
function ToNum($s) {
    return [double]$s
}

# New model results (full reworking of all models)
$params = @("Country[T.France]","Country[T.Germany]","Country[T.Spain]","Age","Commute_Trip","DistCenter_res","DistSubcenter_res","UrbPopDensity_res","UrbBuildDensity_res","IntersecDensity_res","street_length_res","LU_UrbFab_res","LU_Comm_res")

$coeffs = @(
    "-1705.18522655227",
    "-2282.28311162104",
    "-1893.17124822081",
    "25.149781191345699",
    "72.580146680178601",
    "110.14352333005201",
    "56.475632611134102",
    "-3.7927982496068098",
    "-3.1071792806060299",
    "-18.900302625126301",
    "0.797017509519103",
    "-8.9626056031429702",
    "-23.8638599404159"
)

$pvals = @(
    "6.0613959985259301E-8",
    "3.2397184248799E-13",
    "2.66696053516251E-9",
    "3.4383442632439999E-4",
    "1.0663717221661301E-25",
    "1.8619975495636702E-33",
    "3.0011266223089E-3",
    "1.9426894081635601E-7",
    "0.44502679393836703",
    "8.8655069759651101E-13",
    "0.49197492441859902",
    "5.8420730859429996E-4",
    "2.0488828265827199E-11"
)

# Rows 2..14 use the same order/labels on both sheets
$startRow = 2

# ---------------------------------------------------------------------
# Sheet "All_model_short": update B/C columns (coefficient / p), leaving
# existing text labels + header row untouched. Row 10 (UrbBuildDensity_res)
# previously had no coefficient value; now it does.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("All_model_short")

for ($i = 0; $i -lt $params.Count; $i++) {
    $r = $startRow + $i
    $ws1.Cells.Item($r, 2).Value = (ToNum $coeffs[$i])
    $ws1.Cells.Item($r, 3).Value = (ToNum $pvals[$i])
}

# Drop the now unused conditional formatting rule (p between 0.05 and 0.1)
$ws1.Range("C2:C14").FormatConditions.Delete()

# Remove the bespoke column-A width customisation
$ws1.Columns.Item(1).ColumnWidth = $ws1.StandardWidth

# Update the selection shown when the sheet is opened (whole table selected,
# active cell on the last row)
$ws1.Activate()
$ws1.Range("A1:C14").Select()

# ---------------------------------------------------------------------
# Sheet "table" -> rename to "All_model_short_table" and rebuild its
# contents to mirror "All_model_short" (plain A:C table, same labels,
# same values/number formats as the source sheet - this sheet is a
# second, differently-styled copy of the results table).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("table")
$ws2.Name = "All_model_short_table"

$ws2.Cells.Clear()

$headers = @("param","coefficient","p")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $ws2.Cells.Item(1, $c + 1).Value = $headers[$c]
}

for ($i = 0; $i -lt $params.Count; $i++) {
    $r = $startRow + $i
    $ws2.Cells.Item($r, 1).Value = $params[$i]
    $ws2.Cells.Item($r, 2).Value = (ToNum $coeffs[$i])
    $ws2.Cells.Item($r, 3).Value = (ToNum $pvals[$i])
}

# Number formats for the "table" sheet: big/whole coefficients get "0"
# (integer-like) formatting for the first three country rows get "0.0",
# most coefficients get "0.0" style too except the big country ones which
# use plain "0"; p-values mostly use "0.000" with one cell "0.0000".
$ws2.Range("B2:B4").NumberFormat = "0"
$ws2.Range("B5:B14").NumberFormat = "0.0"
$ws2.Range("C2:C14").NumberFormat = "0.000"
$ws2.Range("C13").NumberFormat = "0.0000"

$ws2.Activate()
$ws2.Range("C13").Select()
